$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 500.15
$ws.Range("I96").Value = 292.42856
$ws.Range("J96").Value = 984.8333
$ws.Range("K96").Value = 877.28568
$ws.Range("L96").Value = 2954.4999
$ws.Range("M96").Value = 495.71432
$ws.Range("N96").Value = -5700.4999
$ws.Range("H127").Value = 2130.6033
$ws.Range("J127").Value = 2333.9106
$ws.Range("L127").Value = 7001.7318
$ws.Range("N127").Value = -16921.7318
$ws.Range("H132").Value = 8341226
$ws.Range("I132").Value = 10008594
$ws.Range("J132").Value = 4383.2
$ws.Range("K132").Value = 30025782
$ws.Range("L132").Value = 13149.6
$ws.Range("M132").Value = -30023252
$ws.Range("N132").Value = -18209.6
$ws.Range("H137").Value = 1636.0834
$ws.Range("I137").Value = 1282.5172
$ws.Range("J137").Value = 3100.8572
$ws.Range("K137").Value = 3847.5516
$ws.Range("L137").Value = 9302.571599999999
$ws.Range("M137").Value = -1297.5516
$ws.Range("N137").Value = -14402.5716

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1708.8846
$ws.Range("I61").Value = 1634.8096
$ws.Range("K61").Value = 1634.8096
$ws.Range("M61").Value = -1422.8096
$ws.Range("H74").Value = 2578.9688
$ws.Range("I74").Value = 1365.3914
$ws.Range("J74").Value = 5680.3335
$ws.Range("K74").Value = 1365.3914
$ws.Range("L74").Value = 5680.3335
$ws.Range("M74").Value = -491.3914
$ws.Range("N74").Value = -7428.3335
$ws.Range("H77").Value = 2578.9688
$ws.Range("I77").Value = 1365.3914
$ws.Range("J77").Value = 5680.3335
$ws.Range("K77").Value = 6826.957
$ws.Range("L77").Value = 28401.6675
$ws.Range("M77").Value = -2458.957
$ws.Range("N77").Value = -37137.6675
$ws.Range("H98").Value = 16407.334
$ws.Range("J98").Value = 16407.334
$ws.Range("L98").Value = 16407.334
$ws.Range("N98").Value = -22397.334
$ws.Range("H118").Value = 39385
$ws.Range("J118").Value = 39385
$ws.Range("L118").Value = 39385
$ws.Range("N118").Value = -42699
$ws.Range("H132").Value = 1783.9259
$ws.Range("I132").Value = 1418.25
$ws.Range("K132").Value = 4254.75
$ws.Range("M132").Value = -1724.75
$ws.Range("H136").Value = 1708.8846
$ws.Range("I136").Value = 1634.8096
$ws.Range("K136").Value = 4904.4288
$ws.Range("M136").Value = -2354.4288

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 18333.334
$ws.Range("I26").Value = 18333.334
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 18333.334
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -18041.334
$ws.Range("N26").ClearContents()
$ws.Range("H134").Value = 2136.8462
$ws.Range("I134").Value = 2178.4167
$ws.Range("J134").Value = 1638
$ws.Range("K134").Value = 6535.250100000001
$ws.Range("L134").Value = 4914
$ws.Range("M134").Value = -4000.250100000001
$ws.Range("N134").Value = -9984

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 45445.883
$ws.Range("I31").Value = 39999.46
$ws.Range("K31").Value = 39999.46
$ws.Range("M31").Value = -39704.46
$ws.Range("H34").Value = 45445.883
$ws.Range("I34").Value = 39999.46
$ws.Range("K34").Value = 39999.46
$ws.Range("M34").Value = -39797.46
$ws.Range("H58").Value = 6307.1763
$ws.Range("I58").Value = 1521.8948
$ws.Range("J58").Value = 12368.533
$ws.Range("K58").Value = 1521.8948
$ws.Range("L58").Value = 12368.533
$ws.Range("M58").Value = -1318.8948
$ws.Range("N58").Value = -12774.533
$ws.Range("H132").Value = 4713.9614
$ws.Range("I132").Value = 4825.5
$ws.Range("J132").Value = 4535.5
$ws.Range("K132").Value = 14476.5
$ws.Range("L132").Value = 13606.5
$ws.Range("M132").Value = -11946.5
$ws.Range("N132").Value = -18666.5
$ws.Range("H134").Value = 1256.9395
$ws.Range("I134").Value = 1163.4615
$ws.Range("J134").Value = 1604.1428
$ws.Range("K134").Value = 3490.3845
$ws.Range("L134").Value = 4812.428400000001
$ws.Range("M134").Value = -955.3844999999997
$ws.Range("N134").Value = -9882.428400000001
$ws.Range("H136").Value = 6307.1763
$ws.Range("I136").Value = 1521.8948
$ws.Range("J136").Value = 12368.533
$ws.Range("K136").Value = 4565.6844
$ws.Range("L136").Value = 37105.599
$ws.Range("M136").Value = -2015.6844
$ws.Range("N136").Value = -42205.599
$ws.Range("H138").Value = 67092.5
$ws.Range("J138").Value = 67092.5
$ws.Range("L138").Value = 67092.5
$ws.Range("N138").Value = -77372.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 24719.953
$ws.Range("J37").Value = 24719.953
$ws.Range("L37").Value = 74159.859
$ws.Range("N37").Value = -74383.859
$ws.Range("H40").Value = 348.85715
$ws.Range("I40").Value = 393.66666
$ws.Range("J40").Value = 80
$ws.Range("K40").Value = 1574.66664
$ws.Range("L40").Value = 320
$ws.Range("M40").Value = -1505.66664
$ws.Range("N40").Value = -458
$ws.Range("H113").Value = 563.3043
$ws.Range("I113").Value = 467.875
$ws.Range("J113").Value = 614.2
$ws.Range("K113").Value = 1403.625
$ws.Range("L113").Value = 1842.6
$ws.Range("M113").Value = 766.375
$ws.Range("N113").Value = -6182.6
$ws.Range("H120").Value = 94034.45
$ws.Range("I120").Value = 169063.17
$ws.Range("J120").Value = 4000
$ws.Range("K120").Value = 507189.51
$ws.Range("L120").Value = 12000
$ws.Range("M120").Value = -502351.51
$ws.Range("N120").Value = -21676
$ws.Range("H123").Value = 4066.2727
$ws.Range("I123").Value = 2626.6667
$ws.Range("J123").Value = 4606.125
$ws.Range("K123").Value = 7880.000100000001
$ws.Range("L123").Value = 13818.375
$ws.Range("M123").Value = -5430.000100000001
$ws.Range("N123").Value = -18718.375

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 4444.4443
$ws.Range("J92").Value = 4444.4443
$ws.Range("L92").Value = 4444.4443
$ws.Range("N92").Value = -8188.4443
$ws.Range("H122").Value = 1064.4445
$ws.Range("I122").Value = 950
$ws.Range("K122").Value = 2850
$ws.Range("M122").Value = -400
$ws.Range("H132").Value = 2099.9033
$ws.Range("I132").Value = 1537.1052
$ws.Range("J132").Value = 2991
$ws.Range("K132").Value = 4611.3156
$ws.Range("L132").Value = 8973
$ws.Range("M132").Value = -2081.3156
$ws.Range("N132").Value = -14033

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 810.96
$ws.Range("I22").Value = 749.5
$ws.Range("J22").Value = 822.6667
$ws.Range("K22").Value = 749.5
$ws.Range("L22").Value = 822.6667
$ws.Range("M22").Value = -454.5
$ws.Range("N22").Value = -1412.6667
$ws.Range("H27").Value = 810.96
$ws.Range("I27").Value = 749.5
$ws.Range("J27").Value = 822.6667
$ws.Range("K27").Value = 749.5
$ws.Range("L27").Value = 822.6667
$ws.Range("M27").Value = -642.5
$ws.Range("N27").Value = -1036.6667
$ws.Range("H136").Value = 1254.3784
$ws.Range("I136").Value = 1121.7333
$ws.Range("J136").Value = 1822.8572
$ws.Range("K136").Value = 3365.199900000001
$ws.Range("L136").Value = 5468.571599999999
$ws.Range("M136").Value = -815.1999000000005
$ws.Range("N136").Value = -10568.5716

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H132").Value = 2933.8438
$ws.Range("I132").Value = 2888.5715
$ws.Range("J132").Value = 3250.75
$ws.Range("K132").Value = 8665.7145
$ws.Range("L132").Value = 9752.25
$ws.Range("M132").Value = -6135.7145
$ws.Range("N132").Value = -14812.25
$ws.Range("H136").Value = 1090.6666
$ws.Range("I136").Value = 816
$ws.Range("J136").Value = 1640
$ws.Range("K136").Value = 2448
$ws.Range("L136").Value = 4920
$ws.Range("M136").Value = 102
$ws.Range("N136").Value = -10020
